# 1. Rename the original sheet and add the new QuickSort sheet after it.
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "Cambio Residuo"

$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "QuickSort"

# 2. Populate column A of QuickSort with the quickSort timing-metric values.
$quickSortValues = @("2.9999999999999999E-7","6.9999999999999997E-7","7.9999999999999996E-7","1.3E-6","1.9999999999999999E-6","2.2000000000000001E-6","2.6000000000000001E-6","3.1E-6","3.4000000000000001E-6","4.0999999999999997E-6","4.4000000000000002E-6","4.8999999999999997E-6","5.6999999999999996E-6","6.1E-6","6.7000000000000002E-6","6.9E-6","7.7000000000000008E-6","7.7000000000000008E-6","8.1999999999999994E-6","9.0999999999999993E-6","9.3999999999999998E-6","1.01E-5","1.08E-5","1.1600000000000001E-5","1.2099999999999999E-5","1.2099999999999999E-5","1.33E-5","1.3699999999999999E-5","1.38E-5","1.5E-5","1.7499999999999998E-5","1.9300000000000002E-5","1.9400000000000001E-5","2.1299999999999999E-5","1.6799999999999998E-5","1.7600000000000001E-5","1.84E-5","1.8600000000000001E-5","1.8499999999999999E-5","2.0299999999999999E-5","2.1999999999999999E-5","2.0999999999999999E-5","2.19E-5","2.2799999999999999E-5","2.2900000000000001E-5","2.34E-5","2.3499999999999999E-5","2.4600000000000002E-5","2.44E-5","2.5599999999999999E-5") | ForEach-Object { [double]$_ }
for ($i = 0; $i -lt $quickSortValues.Length; $i++) {
    $ws2.Cells.Item($i + 1, 1).Value = $quickSortValues[$i]
}

# 3. Give column A an explicit "General" number format + a slightly wider column,
#    matching the look of the authored workbook.
$ws2.Columns.Item(1).NumberFormat = "General"
$ws2.Columns.Item(1).ColumnWidth = 10.65

# 4. Select the whole column, as the author had it selected when the file was saved.
$ws2.Columns.Item(1).Select()

# 5. Re-point the two existing charts (on "Cambio Residuo") from the old sheet name
#    to the new one, since chart series formulas store the sheet name literally.
for ($i = 1; $i -le $ws1.ChartObjects().Count; $i++) {
    $co = $ws1.ChartObjects($i)
    $ch = $co.Chart
    $ser = $ch.SeriesCollection(1)
    $ser.Formula = $ser.Formula.Replace("Hoja1!", "'Cambio Residuo'!")
}

# 6. Add a new line chart on QuickSort plotting the new column.
$shp = $ws2.Shapes.AddChart2(201, 4)
$chartObj = $shp.Chart
$chartObj.SetSourceData($ws2.Range("A1:A50"))
$newSer = $chartObj.SeriesCollection(1)
$newSer.Formula = "=SERIES(,,QuickSort!`$A`$1:`$A`$50,1)"

# 7. QuickSort is now the active sheet/tab (Worksheets.Add makes the new sheet active,
#    matching the workbook's activeTab + tabSelected changes).
$ws2.Activate()
